$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 3) matching the header/data layout of row 2
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 42605.885844907411

$ws.Range("B3").Value = -24
$ws.Range("C3").Value = 58
$ws.Range("D3").Value = 41
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 17294
$ws.Range("H3").Value = 3626
$ws.Range("I3").Value = 231
$ws.Range("J3").Value = 27
$ws.Range("K3").Value = 19
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = "Named"
